$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy H1's formatting (bold,
# bordered, centered header style) so the new cells share the same style.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells I2 and J2
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
